# Refine club name matching:
#   "Balance Tri Club" -> "Balance Triathlon Club"
#   "South West Sydney Triathlon Club" -> "Sydney South West Triathlon Club"
#
# Scan the Club Name column (column C) on the active worksheet and
# replace the old names with the corrected ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()

    if ($val -eq "Balance Tri Club") {
        $cell.Value = "Balance Triathlon Club"
    }
    elseif ($val -eq "South West Sydney Triathlon Club") {
        $cell.Value = "Sydney South West Triathlon Club"
    }
}
